$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A76").Value = "01-07-2021"
$values = @(24763,22505,7557,3026,902,343,566,323,312,1933,814,312,2597,637,450,126,80,129,292,341,12598,3706,1194,1069,382,589,51,133,8892,1321,456,1156,185,762,123,171,322,108,180,286,4607,691,93,55,302,328,55,148,2,1303,278,307,272)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(76, $i + 2).Value = $values[$i]
}
